{"js": "// Update the worksheet date and the 25 division-problem answers.\n// The date paragraph text changes, and 25 table cells (5 \"fact rows\" of 5\n// columns each, interleaved with blank rows) get new \"a\u00f7b=c, d\" strings.\n// We update run text in place (via insertText/value) so existing run\n// formatting (fonts, size, etc.) is preserved.\n\nconst dateUpdate = { oldText: \"2023-12-10 Sunday\", newText: \"2023-12-11 Monday\" };\n\n// row/col are 0-based indices into Table.getCell(row, col); the table has\n// 20 rows total, but only rows 0, 4, 8, 12, 16 hold problems (the rest are\n// blank answer rows), each with 5 columns.\nconst cellUpdates = [\n  { row: 0, col: 0, oldText: \"57\u00f79=6, 3\", newText: \"77\u00f78=9, 5\" },\n  { row: 0, col: 1, oldText: \"94\u00f74=23, 2\", newText: \"38\u00f76=6, 2\" },\n  { row: 0, col: 2, oldText: \"88\u00f79=9, 7\", newText: \"42\u00f73=14, 0\" },\n  { row: 0, col: 3, oldText: \"45\u00f72=22, 1\", newText: \"33\u00f73=11, 0\" },\n  { row: 0, col: 4, oldText: \"59\u00f72=29, 1\", newText: \"62\u00f75=12, 2\" },\n  { row: 4, col: 0, oldText: \"36\u00f76=6, 0\", newText: \"15\u00f75=3, 0\" },\n  { row: 4, col: 1, oldText: \"75\u00f78=9, 3\", newText: \"30\u00f72=15, 0\" },\n  { row: 4, col: 2, oldText: \"35\u00f79=3, 8\", newText: \"62\u00f74=15, 2\" },\n  { row: 4, col: 3, oldText: \"69\u00f79=7, 6\", newText: \"68\u00f77=9, 5\" },\n  { row: 4, col: 4, oldText: \"84\u00f79=9, 3\", newText: \"85\u00f72=42, 1\" },\n  { row: 8, col: 0, oldText: \"70\u00f77=10, 0\", newText: \"14\u00f74=3, 2\" },\n  { row: 8, col: 1, oldText: \"50\u00f72=25, 0\", newText: \"57\u00f73=19, 0\" },\n  { row: 8, col: 2, oldText: \"96\u00f79=10, 6\", newText: \"15\u00f77=2, 1\" },\n  { row: 8, col: 3, oldText: \"57\u00f79=6, 3\", newText: \"71\u00f77=10, 1\" },\n  { row: 8, col: 4, oldText: \"51\u00f79=5, 6\", newText: \"85\u00f79=9, 4\" },\n  { row: 12, col: 0, oldText: \"31\u00f78=3, 7\", newText: \"93\u00f78=11, 5\" },\n  { row: 12, col: 1, oldText: \"18\u00f77=2, 4\", newText: \"52\u00f73=17, 1\" },\n  { row: 12, col: 2, oldText: \"82\u00f77=11, 5\", newText: \"88\u00f74=22, 0\" },\n  { row: 12, col: 3, oldText: \"81\u00f78=10, 1\", newText: \"96\u00f76=16, 0\" },\n  { row: 12, col: 4, oldText: \"56\u00f78=7, 0\", newText: \"67\u00f73=22, 1\" },\n  { row: 16, col: 0, oldText: \"10\u00f77=1, 3\", newText: \"56\u00f73=18, 2\" },\n  { row: 16, col: 1, oldText: \"11\u00f76=1, 5\", newText: \"19\u00f79=2, 1\" },\n  { row: 16, col: 2, oldText: \"70\u00f76=11, 4\", newText: \"69\u00f78=8, 5\" },\n  { row: 16, col: 3, oldText: \"75\u00f79=8, 3\", newText: \"98\u00f72=49, 0\" },\n  { row: 16, col: 4, oldText: \"66\u00f78=8, 2\", newText: \"37\u00f76=6, 1\" },\n];\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph (first paragraph of the document) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet dateParagraph = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text === dateUpdate.oldText) {\n    dateParagraph = p;\n    break;\n  }\n}\nif (!dateParagraph && paragraphs.items.length > 0) {\n  dateParagraph = paragraphs.items[0];\n}\nif (dateParagraph) {\n  dateParagraph.insertText(dateUpdate.newText, \"Replace\");\n}\n\n// --- 2. Update the table cells with the new division problems ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const update of cellUpdates) {\n  const cell = table.getCell(update.row, update.col);\n  cell.value = update.newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 division-problem answers.\n# The date paragraph text changes, and 25 table cells (5 \"fact rows\" of 5\n# columns each, interleaved with blank rows) get new \"a\u00f7b=c, d\" strings.\n# We overwrite Range.Text in place so existing run formatting\n# (fonts, size, etc.) carried by the paragraph/run is preserved.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph (first paragraph of the document) ---\n$dateOld = \"2023-12-10 Sunday\"\n$dateNew = \"2023-12-11 Monday\"\n\n$p1 = $d.Paragraphs.Item(1)\n$p1.Range.Text = $dateNew\n\n# --- 2. Update the table cells with the new division problems ---\n# Row/Col are 1-based (Word COM Table.Cell(row, col) indexing). The table\n# has 20 rows total, but only rows 1, 5, 9, 13, 17 hold problems (the rest\n# are blank answer rows), each with 5 columns.\n$cellUpdates = @(\n  @{ Row = 1; Col = 1; OldText = \"57\u00f79=6, 3\"; NewText = \"77\u00f78=9, 5\" },\n  @{ Row = 1; Col = 2; OldText = \"94\u00f74=23, 2\"; NewText = \"38\u00f76=6, 2\" },\n  @{ Row = 1; Col = 3; OldText = \"88\u00f79=9, 7\"; NewText = \"42\u00f73=14, 0\" },\n  @{ Row = 1; Col = 4; OldText = \"45\u00f72=22, 1\"; NewText = \"33\u00f73=11, 0\" },\n  @{ Row = 1; Col = 5; OldText = \"59\u00f72=29, 1\"; NewText = \"62\u00f75=12, 2\" },\n  @{ Row = 5; Col = 1; OldText = \"36\u00f76=6, 0\"; NewText = \"15\u00f75=3, 0\" },\n  @{ Row = 5; Col = 2; OldText = \"75\u00f78=9, 3\"; NewText = \"30\u00f72=15, 0\" },\n  @{ Row = 5; Col = 3; OldText = \"35\u00f79=3, 8\"; NewText = \"62\u00f74=15, 2\" },\n  @{ Row = 5; Col = 4; OldText = \"69\u00f79=7, 6\"; NewText = \"68\u00f77=9, 5\" },\n  @{ Row = 5; Col = 5; OldText = \"84\u00f79=9, 3\"; NewText = \"85\u00f72=42, 1\" },\n  @{ Row = 9; Col = 1; OldText = \"70\u00f77=10, 0\"; NewText = \"14\u00f74=3, 2\" },\n  @{ Row = 9; Col = 2; OldText = \"50\u00f72=25, 0\"; NewText = \"57\u00f73=19, 0\" },\n  @{ Row = 9; Col = 3; OldText = \"96\u00f79=10, 6\"; NewText = \"15\u00f77=2, 1\" },\n  @{ Row = 9; Col = 4; OldText = \"57\u00f79=6, 3\"; NewText = \"71\u00f77=10, 1\" },\n  @{ Row = 9; Col = 5; OldText = \"51\u00f79=5, 6\"; NewText = \"85\u00f79=9, 4\" },\n  @{ Row = 13; Col = 1; OldText = \"31\u00f78=3, 7\"; NewText = \"93\u00f78=11, 5\" },\n  @{ Row = 13; Col = 2; OldText = \"18\u00f77=2, 4\"; NewText = \"52\u00f73=17, 1\" },\n  @{ Row = 13; Col = 3; OldText = \"82\u00f77=11, 5\"; NewText = \"88\u00f74=22, 0\" },\n  @{ Row = 13; Col = 4; OldText = \"81\u00f78=10, 1\"; NewText = \"96\u00f76=16, 0\" },\n  @{ Row = 13; Col = 5; OldText = \"56\u00f78=7, 0\"; NewText = \"67\u00f73=22, 1\" },\n  @{ Row = 17; Col = 1; OldText = \"10\u00f77=1, 3\"; NewText = \"56\u00f73=18, 2\" },\n  @{ Row = 17; Col = 2; OldText = \"11\u00f76=1, 5\"; NewText = \"19\u00f79=2, 1\" },\n  @{ Row = 17; Col = 3; OldText = \"70\u00f76=11, 4\"; NewText = \"69\u00f78=8, 5\" },\n  @{ Row = 17; Col = 4; OldText = \"75\u00f79=8, 3\"; NewText = \"98\u00f72=49, 0\" },\n  @{ Row = 17; Col = 5; OldText = \"66\u00f78=8, 2\"; NewText = \"37\u00f76=6, 1\" }\n)\n\n$t = $d.Tables.Item(1)\n\nforeach ($u in $cellUpdates) {\n  $cell = $t.Cell($u.Row, $u.Col)\n  $cell.Range.Text = $u.NewText\n}\n"}
